$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 379 (shifts old rows 379-427 down to 380-428)
$ws.Rows.Item(379).Insert()

# Column A: timestamp (numeric)
$ws.Cells.Item(379, 1).Value = 1574899200

# Column B: date string "2019-11-28" - force text so it isn't parsed as a date
$ws.Cells.Item(379, 2).NumberFormat = "@"
$ws.Cells.Item(379, 2).Value = "2019-11-28"
$ws.Cells.Item(379, 2).Style = "Normal"

# Column C: id "03005" - force text so the leading zero is preserved
$ws.Cells.Item(379, 3).NumberFormat = "@"
$ws.Cells.Item(379, 3).Value = "03005"
$ws.Cells.Item(379, 3).Style = "Normal"

# Column D: name (plain text)
$ws.Cells.Item(379, 4).Value = "MHCARE"

# Columns E-H: open/high/low/close
$ws.Cells.Item(379, 5).Value = 1.02
$ws.Cells.Item(379, 6).Value = 1.02
$ws.Cells.Item(379, 7).Value = 1.02
$ws.Cells.Item(379, 8).Value = 1.02

# Column I: vol
$ws.Cells.Item(379, 9).Value = 500
